$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9127689003944397
$ws.Range("B1").Value = 3.107428550720215
$ws.Range("C1").Value = 2.818970441818237
$ws.Range("D1").Value = 1.60148811340332
$ws.Range("E1").Value = 1.230545878410339
